# Populate row 1 of Sheet1 with the new demo values.
# Cell-entry order below mirrors the order the strings first appear in the
# saved shared-strings table (A1 is numeric so it doesn't affect that table;
# K1's "," is entered right after B1 so it lands at shared-string index 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4
$ws.Range("B1").Value = "4``"
$ws.Range("K1").Value = ","
$ws.Range("C1").Value = "m"
$ws.Range("D1").Value = "m"
$ws.Range("E1").Value = "l"
$ws.Range("F1").Value = "l"
$ws.Range("G1").Value = "l"
$ws.Range("H1").Value = "l"
$ws.Range("I1").Value = "l"
$ws.Range("J1").Value = "l"
$ws.Range("L1").Value = "k"
$ws.Range("M1").Value = "n"
$ws.Range("N1").Value = "knk"
$ws.Range("P1").Value = "lk"
$ws.Range("Q1").Value = "ml"
$ws.Range("S1").Value = "lm"

# Leave the selection on S1, matching the saved workbook's cursor position.
$ws.Range("S1").Select() | Out-Null
